# Fix Training Data Issue (#48)
# The "Date" column (BF) on the sheet held the literal text "5-2-2011-12"
# for every data row. NBA stats for the date in question were actually
# reported a day off, so the column needs to read "2012-05-02" instead.
#
# The values must stay plain text (they are not real Excel dates), so the
# number format is forced to Text ("@") before the write and the cell
# style is reset back to Normal afterwards so no stray formatting is left
# behind on the cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$range = $ws.Range("BF2:BF31")
$range.NumberFormat = "@"
$range.Value = "2012-05-02"
$range.Style = "Normal"
